$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (shifts existing rows 15-64 down to 16-65)
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new weekly record
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 45063
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100112030
$ws.Range("G15").Value = "Poroto granado"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 30000
$ws.Range("L15").Value = 32000
$ws.Range("M15").Value = 31000
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 1240
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
